$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Edit 1: "Team" -> "Team ID" (append a new run " ID" after "Team") ---
$teamCell = $t.Cell(2, 1)
$teamRange = $teamCell.Range
$insPoint = $teamRange.End - 1   # just after "Team", before the cell mark
$ins = $d.Range($insPoint, $insPoint)
$ins.InsertAfter(" ID")

# Touch Bold (true then false) on the freshly inserted text so it stays a
# distinct run (matching Word's own run-splitting behaviour) while ending
# up with the exact same, non-bold formatting as the "Team" run.
$newRun = $d.Range($insPoint, $insPoint + 3)
$newRun.Bold = 1
$newRun.Bold = 0

# --- Edit 2: "AS PS VS VV" -> "178047" (and drop the bold formatting) ---
$d.Content.Find.Execute("AS PS VS VV", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "178047", 2)

$teamIdCell = $t.Cell(2, 2)
$valueRange = $d.Range($teamIdCell.Range.Start, $teamIdCell.Range.End - 1)
$valueRange.Font.Bold = 0
